$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3361
$ws.Range("E2").Value = 282
$ws.Range("F2").Value = 282
$ws.Range("G2").Value = 249
$ws.Range("H2").Value = 214
$ws.Range("I2").Value = 214
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 4744
$ws.Range("L2").Value = 2611
$ws.Range("M2").Value = 2132
$ws.Range("N2").Value = 2129
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 565
$ws.Range("Q2").Value = 324
$ws.Range("R2").Value = -456
$ws.Range("S2").Value = 160
$ws.Range("T2").Value = 439
$ws.Range("U2").Value = -115
$ws.Range("V2").Value = 1797
$ws.Range("W2").Value = 8.390000000000001
$ws.Range("X2").Value = 6.37
$ws.Range("Y2").Value = 10.42
$ws.Range("Z2").Value = 4.75
$ws.Range("AA2").Value = 122.45
$ws.Range("AB2").Value = 270.08
$ws.Range("AC2").Value = 1897
$ws.Range("AD2").Value = 20.54
$ws.Range("AE2").Value = 19203
$ws.Range("AF2").Value = 2.03
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 1.54
$ws.Range("AI2").Value = 31.06
$ws.Range("AJ2").Value = 11295195

# Row 3
$ws.Range("D3").Value = 3680
$ws.Range("E3").Value = 490
$ws.Range("F3").Value = 490
$ws.Range("G3").Value = 426
$ws.Range("H3").Value = 343
$ws.Range("I3").Value = 341
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 5210
$ws.Range("L3").Value = 2743
$ws.Range("M3").Value = 2468
$ws.Range("N3").Value = 2420
$ws.Range("O3").Value = 47
$ws.Range("P3").Value = 565
$ws.Range("Q3").Value = 473
$ws.Range("R3").Value = -616
$ws.Range("S3").Value = 98
$ws.Range("T3").Value = 586
$ws.Range("U3").Value = -113
$ws.Range("V3").Value = 1981
$ws.Range("W3").Value = 13.31
$ws.Range("X3").Value = 9.32
$ws.Range("Y3").Value = 15.01
$ws.Range("Z3").Value = 6.89
$ws.Range("AA3").Value = 111.14
$ws.Range("AB3").Value = 316.74
$ws.Range("AC3").Value = 3023
$ws.Range("AD3").Value = 20.54
$ws.Range("AE3").Value = 21827
$ws.Range("AF3").Value = 2.85
$ws.Range("AG3").Value = 800
$ws.Range("AH3").Value = 1.29
$ws.Range("AI3").Value = 25.98
$ws.Range("AJ3").Value = 11295195

# Row 4
$ws.Range("D4").Value = 4604
$ws.Range("E4").Value = 821
$ws.Range("F4").Value = 821
$ws.Range("G4").Value = 774
$ws.Range("H4").Value = 589
$ws.Range("I4").Value = 558
$ws.Range("J4").Value = 31
$ws.Range("K4").Value = 6997
$ws.Range("L4").Value = 3586
$ws.Range("M4").Value = 3411
$ws.Range("N4").Value = 2824
$ws.Range("O4").Value = 588
$ws.Range("P4").Value = 565
$ws.Range("Q4").Value = 969
$ws.Range("R4").Value = -1468
$ws.Range("S4").Value = 626
$ws.Range("T4").Value = 384
$ws.Range("U4").Value = 585
$ws.Range("V4").Value = 2502
$ws.Range("W4").Value = 17.84
$ws.Range("X4").Value = 12.8
$ws.Range("Y4").Value = 21.29
$ws.Range("Z4").Value = 9.65
$ws.Range("AA4").Value = 105.12
$ws.Range("AB4").Value = 398.76
$ws.Range("AC4").Value = 4942
$ws.Range("AD4").Value = 16.86
$ws.Range("AE4").Value = 25465
$ws.Range("AF4").Value = 3.27
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1.2
$ws.Range("AI4").Value = 19.86
$ws.Range("AJ4").Value = 11295195

# Row 5
$ws.Range("D5").Value = 5216
$ws.Range("E5").Value = 792
$ws.Range("F5").Value = 792
$ws.Range("G5").Value = 747
$ws.Range("H5").Value = 581
$ws.Range("I5").Value = 543
$ws.Range("J5").Value = 38
$ws.Range("K5").Value = 7233
$ws.Range("L5").Value = 3325
$ws.Range("M5").Value = 3909
$ws.Range("N5").Value = 3207
$ws.Range("O5").Value = 701
$ws.Range("P5").Value = 565
$ws.Range("Q5").Value = 874
$ws.Range("R5").Value = -329
$ws.Range("S5").Value = -295
$ws.Range("T5").Value = 310
$ws.Range("U5").Value = 564
$ws.Range("V5").Value = 2185
$ws.Range("W5").Value = 15.18
$ws.Range("X5").Value = 11.13
$ws.Range("Y5").Value = 17.99
$ws.Range("Z5").Value = 8.16
$ws.Range("AA5").Value = 85.06
$ws.Range("AB5").Value = 473.38
$ws.Range("AC5").Value = 4804
$ws.Range("AD5").Value = 15.09
$ws.Range("AE5").Value = 28926
$ws.Range("AF5").Value = 2.51
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1.38
$ws.Range("AI5").Value = 20.44
$ws.Range("AJ5").Value = 11295195

# Row 6
$ws.Range("D6").Value = 5819
$ws.Range("E6").Value = 936
$ws.Range("F6").Value = 936
$ws.Range("G6").Value = 925
$ws.Range("H6").Value = 678
$ws.Range("I6").Value = 705
$ws.Range("K6").Value = 7917
$ws.Range("L6").Value = 3445
$ws.Range("M6").Value = 4472
$ws.Range("N6").Value = 3798
$ws.Range("P6").Value = 565
$ws.Range("Q6").Value = 1006
$ws.Range("R6").Value = -1215
$ws.Range("S6").Value = -146
$ws.Range("T6").Value = 1182
$ws.Range("U6").Value = -177
$ws.Range("V6").Value = 2168
$ws.Range("W6").Value = 16.08
$ws.Range("X6").Value = 11.64
$ws.Range("Y6").Value = 20.12
$ws.Range("Z6").Value = 8.94
$ws.Range("AA6").Value = 77.03
$ws.Range("AB6").Value = 590.15
$ws.Range("AC6").Value = 6238
$ws.Range("AD6").Value = 12.42
$ws.Range("AE6").Value = 34251
$ws.Range("AF6").Value = 2.26
$ws.Range("AG6").Value = 1200
$ws.Range("AH6").Value = 1.55
$ws.Range("AI6").Value = 18.88
$ws.Range("AJ6").Value = 11295195

# Row 7
$ws.Range("D7").Value = 5979
$ws.Range("E7").Value = 1144
$ws.Range("G7").Value = 1158
$ws.Range("H7").Value = 890
$ws.Range("I7").Value = 885
$ws.Range("K7").Value = 8924
$ws.Range("L7").Value = 3699
$ws.Range("M7").Value = 5225
$ws.Range("N7").Value = 4543
$ws.Range("P7").Value = 561
$ws.Range("Q7").Value = 1129
$ws.Range("R7").Value = -649
$ws.Range("S7").Value = -18
$ws.Range("T7").Value = 655
$ws.Range("U7").Value = 486
$ws.Range("W7").Value = 19.13
$ws.Range("X7").Value = 14.89
$ws.Range("Y7").Value = 21.23
$ws.Range("Z7").Value = 10.57
$ws.Range("AA7").Value = 70.8
$ws.Range("AC7").Value = 7839
$ws.Range("AD7").Value = 13.91
$ws.Range("AE7").Value = 40973
$ws.Range("AF7").Value = 2.66
$ws.Range("AG7").Value = 1200
$ws.Range("AH7").Value = 1.1
$ws.Range("AI7").Value = 15.31

# Row 8
$ws.Range("D8").Value = 6678
$ws.Range("E8").Value = 1351
$ws.Range("G8").Value = 1376
$ws.Range("H8").Value = 1043
$ws.Range("I8").Value = 1036
$ws.Range("K8").Value = 9869
$ws.Range("L8").Value = 3738
$ws.Range("M8").Value = 6131
$ws.Range("N8").Value = 5442
$ws.Range("P8").Value = 561
$ws.Range("Q8").Value = 1445
$ws.Range("R8").Value = -722
$ws.Range("S8").Value = -220
$ws.Range("T8").Value = 688
$ws.Range("U8").Value = 784
$ws.Range("W8").Value = 20.23
$ws.Range("X8").Value = 15.61
$ws.Range("Y8").Value = 20.74
$ws.Range("Z8").Value = 11.1
$ws.Range("AA8").Value = 60.97
$ws.Range("AC8").Value = 9169
$ws.Range("AD8").Value = 11.89
$ws.Range("AE8").Value = 49075
$ws.Range("AF8").Value = 2.22
$ws.Range("AG8").Value = 1225
$ws.Range("AH8").Value = 1.12
$ws.Range("AI8").Value = 13.36

# Row 9
$ws.Range("D9").Value = 7443
$ws.Range("E9").Value = 1609
$ws.Range("G9").Value = 1651
$ws.Range("H9").Value = 1248
$ws.Range("I9").Value = 1243
$ws.Range("K9").Value = 11083
$ws.Range("L9").Value = 3842
$ws.Range("M9").Value = 7242
$ws.Range("N9").Value = 6547
$ws.Range("P9").Value = 561
$ws.Range("Q9").Value = 1604
$ws.Range("R9").Value = -729
$ws.Range("S9").Value = -135
$ws.Range("T9").Value = 685
$ws.Range("U9").Value = 1099
$ws.Range("W9").Value = 21.61
$ws.Range("X9").Value = 16.77
$ws.Range("Y9").Value = 20.73
$ws.Range("Z9").Value = 11.91
$ws.Range("AA9").Value = 53.05
$ws.Range("AC9").Value = 11004
$ws.Range("AD9").Value = 9.91
$ws.Range("AE9").Value = 59041
$ws.Range("AF9").Value = 1.85
$ws.Range("AG9").Value = 1225
$ws.Range("AH9").Value = 1.12
$ws.Range("AI9").Value = 11.13
